$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (row 1)
$ws.Range("A1").Value = "Datos actualizados a 23 de Mayo de 2020 a las 19:05"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1654449
$ws.Range("C4").Value = 9355
$ws.Range("E4").Value = 1152170
$ws.Range("G4").Value = 349
$ws.Range("H4").Value = 97996

# Row 5 - Brasil
$ws.Range("B5").Value = 340837
$ws.Range("C5").Value = 9947
$ws.Range("E5").Value = 183729
$ws.Range("G5").Value = 630
$ws.Range("H5").Value = 21678

# Row 12 - Turquia
$ws.Range("B12").Value = 155686
$ws.Range("C12").Value = 1186
$ws.Range("D12").Value = 117602
$ws.Range("E12").Value = 33776
$ws.Range("G12").Value = 32
$ws.Range("H12").Value = 4308

# Row 14 - India
$ws.Range("B14").Value = 131103
$ws.Range("C14").Value = 6309
$ws.Range("D14").Value = 54298
$ws.Range("E14").Value = 72940
$ws.Range("G14").Value = 139
$ws.Range("H14").Value = 3865

# Row 35 - Polonia
$ws.Range("B35").Value = 20931
$ws.Range("C35").Value = 312
$ws.Range("E35").Value = 10961
$ws.Range("G35").Value = 11
$ws.Range("H35").Value = 993

# Row 40 - Rumania
$ws.Range("E40").Value = 5497
$ws.Range("G40").Value = 7
$ws.Range("H40").Value = 1173

# Row 54 - Barein
$ws.Range("E54").Value = 4299
$ws.Range("G54").Value = 1
$ws.Range("H54").Value = 13

# Row 101 - Maldivas
$ws.Range("B101").Value = 1313
$ws.Range("C101").Value = 39
$ws.Range("E101").Value = 1200

# Row 170 - Bahamas
$ws.Range("D170").Value = 45
$ws.Range("E170").Value = 41
